# Applies the "riska.xlsx" update: refreshed repayment figures for several
# collectors plus a sheet rename (report re-run, "(4)" -> "(5)").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as TEXT (matches the source data, which stores
# formatted numbers like "620,427.00" as plain text / shared strings rather
# than numeric cells) without leaving a lasting number-format style behind.
function Set-TextValue($rangeRef, $text) {
    $rng = $ws.Range($rangeRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.ClearFormats()
}

# Rename the sheet to reflect the newer export revision.
$ws.Name = "repayment_20250924_20250924 (5)"

# Row 3 - Erlangga Hutama
Set-TextValue "E3" "620,427.00"
$ws.Range("D3").Value = 2
Set-TextValue "G3" "0.44"
$ws.Range("H3").Value = 319

# Row 4 - Ridhoi Berkat Zebua
Set-TextValue "E4" "3,833,320.00"
$ws.Range("D4").Value = 7
Set-TextValue "G4" "2.25"
$ws.Range("H4").Value = 865

# Row 5 - Fadilah Damayanti
$ws.Range("H5").Value = 412

# Row 6 - Adistira Winditya P
$ws.Range("H6").Value = 1.362

# Row 7 - Yandi Nugraha
Set-TextValue "E7" "3,521,531.00"
$ws.Range("D7").Value = 11
Set-TextValue "G7" "2.69"
$ws.Range("H7").Value = 866

# Row 8 - Annisa Putri Restu
$ws.Range("H8").Value = 1.574

# Row 9 - Nuraini
$ws.Range("H9").Value = 1.073
$ws.Range("J9").Value = 1
Set-TextValue "K9" "0.73"
Set-TextValue "L9" "5.88"

# Row 10 - Romli
$ws.Range("H10").Value = 879

# Row 11 - Riska Nurlita
$ws.Range("H11").Value = 532

# Row 12 - Debora Retima Sihombing
$ws.Range("H12").Value = 1.25

# Row 13 - Erick Ervan Dewanggga
$ws.Range("H13").Value = 647

# Row 14 - Aldi Taufik
$ws.Range("H14").Value = 457

# Row 15 - Nur Halim
Set-TextValue "E15" "823,164.00"
$ws.Range("D15").Value = 3
Set-TextValue "G15" "0.59"

# Row 16 - Axl Wicaksono
$ws.Range("H16").Value = 2.533

# Row 17 - Sucika Wardani
$ws.Range("H17").Value = 443

# Row 18 - Wasti Feronika Sihombing
$ws.Range("H18").Value = 814
